$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections -------------------------------------------------
# E12: trailing space added to "SDA"
$ws.Range("E12").Value = "SDA "

# C21: relabeled from "LCD serializer" to "SCD serializer"
$ws.Range("C21").Value = "SCD serializer"

# C7: DHT11 sensor label corrected/expanded
$ws.Range("C7").Value = "DHT11 o DHT Temperature, Humidity and heatindex"

# E21: pin number corrected from 20 to 21
$ws.Range("E21").Value = 21

# --- New cell F17: DHT library type comment ---------------------------
$ws.Range("F17").Value = "HTTYPE DHT21   // AM2301 "

$f = $ws.Range("F17").Font
$f.Name = "Courier New"
$f.Size = 14
$f.Color = 39219
$f.Family = 1

# Row 17 gets a taller custom row height to fit the new font
$ws.Rows(17).RowHeight = 19

# --- Selection as left by the editor -----------------------------------
[void]$ws.Range("H23").Select()
